$d = $word.ActiveDocument

# Find and remove the "pander(table_forecasts)" SourceCode paragraph that
# sits just after the forecast-intro paragraph and just before the table.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*pander(table_forecasts)*") {
        $p.Range.Delete()
        break
    }
}
